$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "list" column header + values for the existing rows
$ws.Range("C1").Value = "list"
$ws.Range("C2").Value = "ruby-tips"
$ws.Range("C3").Value = "vim-tips"
$ws.Range("C4").Value = "ruby-tips"
$ws.Range("C5").Value = "python-tips"

# New rows for scott and alice
$ws.Range("A6").Value = "scott"
$ws.Range("B6").Value = "scott@scott.com"
$ws.Range("C6").Value = "python-tips"
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:scott@scott.com")

$ws.Range("A7").Value = "alice"
$ws.Range("B7").Value = "alice@example.com"
$ws.Range("C7").Value = "no-tips"
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:alice@example.com")

# Match the existing hyperlink look/font for the two newly-added link cells
$ws.Range("B6").Font.Underline = $ws.Range("B2").Font.Underline
$ws.Range("B6").Font.ThemeColor = $ws.Range("B2").Font.ThemeColor
$ws.Range("B7").Font.Underline = $ws.Range("B2").Font.Underline
$ws.Range("B7").Font.ThemeColor = $ws.Range("B2").Font.ThemeColor

# Header row is no longer bold now that the sheet has a proper header style
$ws.Range("A1:C1").Font.Bold = $false

# Column C width (Google-Sheets-style column)
$ws.Columns.Item(3).ColumnWidth = 15.83

# Page setup
$ws.PageSetup.Orientation = 1

# Leave selection where the import script left off
$ws.Range("F4").Select() | Out-Null
